$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 38.04655766666666
$ws.Range("H2").Value = 114.139673
$ws.Range("I2").Value = 0.8090698722086991
$ws.Range("J2").Value = 0.8090698722086992
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 7.818077
$ws.Range("N2").Value = 23.454231
$ws.Range("O2").Value = 0.1819018824829088
$ws.Range("P2").Value = 0.1819018824829088
$ws.Range("Q2").Value = 297.4509174229403
$ws.Range("R2").Value = 2677.058256806463
$ws.Range("S2").Value = 0.1471713328149688
$ws.Range("T2").Value = 0.1471713328149689

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 38.04655766666666
$ws.Range("H3").Value = 114.139673
$ws.Range("I3").Value = 0.8090698722086991
$ws.Range("J3").Value = 0.8090698722086992
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 22.011801
$ws.Range("N3").Value = 66.035403
$ws.Range("O3").Value = 0.5121448712693895
$ws.Range("P3").Value = 0.5121448712693895
$ws.Range("Q3").Value = 837.473256093691
$ws.Range("R3").Value = 7537.259304843219
$ws.Range("S3").Value = 0.4143609855502656
$ws.Range("T3").Value = 0.4143609855502656

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 38.04655766666666
$ws.Range("H4").Value = 114.139673
$ws.Range("I4").Value = 0.8090698722086991
$ws.Range("J4").Value = 0.8090698722086992
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.594576333333333
$ws.Range("N4").Value = 13.783729
$ws.Range("O4").Value = 0.1069012346955337
$ws.Range("P4").Value = 0.1069012346955337
$ws.Range("Q4").Value = 174.8078134200686
$ws.Range("R4").Value = 1573.270320780617
$ws.Range("S4").Value = 0.08649056829406758
$ws.Range("T4").Value = 0.08649056829406762

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 38.04655766666666
$ws.Range("H5").Value = 114.139673
$ws.Range("I5").Value = 0.8090698722086991
$ws.Range("J5").Value = 0.8090698722086992
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 8.555183333333334
$ws.Range("N5").Value = 25.66555
$ws.Range("O5").Value = 0.1990520115521681
$ws.Range("P5").Value = 0.1990520115521681
$ws.Range("Q5").Value = 325.4952760405722
$ws.Range("R5").Value = 2929.45748436515
$ws.Range("S5").Value = 0.1610469855493971
$ws.Range("T5").Value = 0.1610469855493972

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.617245333333334
$ws.Range("H6").Value = 4.851736000000001
$ws.Range("I6").Value = 0.03439113957782537
$ws.Range("J6").Value = 0.03439113957782537
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 7.818077
$ws.Range("N6").Value = 23.454231
$ws.Range("O6").Value = 0.1819018824829088
$ws.Range("P6").Value = 0.1819018824829088
$ws.Range("Q6").Value = 12.64374854389067
$ws.Range("R6").Value = 113.793736895016
$ws.Range("S6").Value = 0.006255813029938904
$ws.Range("T6").Value = 0.006255813029938905

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.617245333333334
$ws.Range("H7").Value = 4.851736000000001
$ws.Range("I7").Value = 0.03439113957782537
$ws.Range("J7").Value = 0.03439113957782537
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 22.011801
$ws.Range("N7").Value = 66.035403
$ws.Range("O7").Value = 0.5121448712693895
$ws.Range("P7").Value = 0.5121448712693895
$ws.Range("Q7").Value = 35.59848244551201
$ws.Range("R7").Value = 320.3863420096081
$ws.Range("S7").Value = 0.01761324575189298
$ws.Range("T7").Value = 0.01761324575189298

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.617245333333334
$ws.Range("H8").Value = 4.851736000000001
$ws.Range("I8").Value = 0.03439113957782537
$ws.Range("J8").Value = 0.03439113957782537
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 4.594576333333333
$ws.Range("N8").Value = 13.783729
$ws.Range("O8").Value = 0.1069012346955337
$ws.Range("P8").Value = 0.1069012346955337
$ws.Range("Q8").Value = 7.430557133727112
$ws.Range("R8").Value = 66.87501420354401
$ws.Range("S8").Value = 0.003676455283455967
$ws.Range("T8").Value = 0.003676455283455968

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.617245333333334
$ws.Range("H9").Value = 4.851736000000001
$ws.Range("I9").Value = 0.03439113957782537
$ws.Range("J9").Value = 0.03439113957782537
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 8.555183333333334
$ws.Range("N9").Value = 25.66555
$ws.Range("O9").Value = 0.1990520115521681
$ws.Range("P9").Value = 0.1990520115521681
$ws.Range("Q9").Value = 13.83583032164445
$ws.Range("R9").Value = 124.5224728948
$ws.Range("S9").Value = 0.006845625512537522
$ws.Range("T9").Value = 0.006845625512537523

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 7.361255
$ws.Range("H10").Value = 22.083765
$ws.Range("I10").Value = 0.1565389882134754
$ws.Range("J10").Value = 0.1565389882134754
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 7.818077
$ws.Range("N10").Value = 23.454231
$ws.Range("O10").Value = 0.1819018824829088
$ws.Range("P10").Value = 0.1819018824829088
$ws.Range("Q10").Value = 57.55085840663499
$ws.Range("R10").Value = 517.957725659715
$ws.Range("S10").Value = 0.02847473663800106
$ws.Range("T10").Value = 0.02847473663800106

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 7.361255
$ws.Range("H11").Value = 22.083765
$ws.Range("I11").Value = 0.1565389882134754
$ws.Range("J11").Value = 0.1565389882134754
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 22.011801
$ws.Range("N11").Value = 66.035403
$ws.Range("O11").Value = 0.5121448712693895
$ws.Range("P11").Value = 0.5121448712693895
$ws.Range("Q11").Value = 162.034480170255
$ws.Range("R11").Value = 1458.310321532295
$ws.Range("S11").Value = 0.08017063996723085
$ws.Range("T11").Value = 0.08017063996723085

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 7.361255
$ws.Range("H12").Value = 22.083765
$ws.Range("I12").Value = 0.1565389882134754
$ws.Range("J12").Value = 0.1565389882134754
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 4.594576333333333
$ws.Range("N12").Value = 13.783729
$ws.Range("O12").Value = 0.1069012346955337
$ws.Range("P12").Value = 0.1069012346955337
$ws.Range("Q12").Value = 33.82184800663167
$ws.Range("R12").Value = 304.396632059685
$ws.Range("S12").Value = 0.01673421111801012
$ws.Range("T12").Value = 0.01673421111801012

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 7.361255
$ws.Range("H13").Value = 22.083765
$ws.Range("I13").Value = 0.1565389882134754
$ws.Range("J13").Value = 0.1565389882134754
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 8.555183333333334
$ws.Range("N13").Value = 25.66555
$ws.Range("O13").Value = 0.1990520115521681
$ws.Range("P13").Value = 0.1990520115521681
$ws.Range("Q13").Value = 62.97688608841667
$ws.Range("R13").Value = 566.79197479575
$ws.Range("S13").Value = 0.03115940049023342
$ws.Range("T13").Value = 0.03115940049023343

